$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Q3/R3: round the coordinate values down to whole numbers (drop the
# fractional part that was previously stored).
$ws.Range("Q3").Value = 555846
$ws.Range("R3").Value = 6952042

# Z3/AB3 ("Starttid"/"Sluttid" = 00:00 time-of-day placeholders) are removed
# entirely for this row, while AA3 ("Slutdatum") is left untouched.
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
